$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "1.028").
# Excel would otherwise auto-convert such text into a real number, which
# silently drops meaningful trailing zeros (e.g. "1.150" -> 1.15). Force
# those specific cells to Text format first so the literal string is kept.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "27.785.38"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "1.863.31"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "323.28"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "1.028"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "0.4403"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").Value = "0.3812"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "0.07443"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "21.68"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "1.880.05"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "5.561"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "6.756"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "0.07181"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "85.62"
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").Value = "1.035"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "0.000009122"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "1.028"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "15.56"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "27.862.91"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "5.309"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "11.29"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "2.108.28"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "2.031"
$ws.Range("E25").Value = "  +7.05%  "
$ws.Range("D26").Value = "158.41"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "18.81"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Value = "5.406"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").Value = "2.005"
$ws.Range("E29").Value = "  +4.53%  "
$ws.Range("D30").Value = "117.95"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "0.09029"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "0.7856"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "1.225"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").Value = "3.023"
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("D35").Value = "4.595"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").Value = "1.029"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "1.150"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "0.01986"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "0.05326"
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("D40").Value = "2.871"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("D41").Value = "0.5222"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "0.1692"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").Value = "6.920"
$ws.Range("E43").Value = "  +5.44%  "
$ws.Range("D44").Value = "8.912"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("D45").Value = "110.93"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "10.74"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "0.06610"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("D48").Value = "1.029"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "1.719"
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("D50").Value = "0.4745"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("E51").Value = "  +1.41%  "

Write-Host "Updated cryptos list"